$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    6  = -3
    10 = -5
    12 = -2
    14 = -4
    17 = 7
    18 = -1
    20 = -6
    23 = -3
    24 = -5
    26 = -4
    33 = -5
    37 = -9
    45 = -2
    48 = -2
    50 = -7
    51 = -5
    53 = -4
    55 = -1
    61 = 3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
